# Apply updated GHI values for RAAL Production model (2024-02-09)
$wb = $excel.ActiveWorkbook

# ---- Daily sheet ----
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 2699.47
$daily.Range("H2").Value = 5949.01
$daily.Range("I2").Value = 693.33
$daily.Range("J2").Value = 679.4400000000001
$daily.Range("L2").Value = 679.4400000000001

# ---- Hourly sheet ----
$hourly = $wb.Worksheets.Item("Hourly")

$hourly.Range("I9").Value = 26.7
$hourly.Range("K9").Value = 1.07
$hourly.Range("M9").Value = 1.07

$hourly.Range("I10").Value = 387.51
$hourly.Range("K10").Value = 26.26
$hourly.Range("M10").Value = 26.26

$hourly.Range("H11").Value = 224.68
$hourly.Range("I11").Value = 612.09

$hourly.Range("I12").Value = 716.12
$hourly.Range("K12").Value = 85.75
$hourly.Range("M12").Value = 85.75

$hourly.Range("H13").Value = 422.63
$hourly.Range("I13").Value = 767.11
$hourly.Range("J13").Value = 90.67
$hourly.Range("K13").Value = 105.66
$hourly.Range("M13").Value = 105.66

$hourly.Range("H14").Value = 453.54
$hourly.Range("I14").Value = 784.25
$hourly.Range("J14").Value = 93.38
$hourly.Range("K14").Value = 113.39
$hourly.Range("M14").Value = 113.39

$hourly.Range("H15").Value = 432.2
$hourly.Range("I15").Value = 772.61
$hourly.Range("K15").Value = 108.05
$hourly.Range("M15").Value = 108.05

$hourly.Range("I16").Value = 728.77
$hourly.Range("J16").Value = 84.81
$hourly.Range("K16").Value = 90.26000000000001
$hourly.Range("M16").Value = 90.26000000000001

$hourly.Range("H17").Value = 248.77
$hourly.Range("I17").Value = 637.26
$hourly.Range("K17").Value = 62.19
$hourly.Range("M17").Value = 62.19

$hourly.Range("I18").Value = 446.01
